# Update EPEX Spot prices workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add column BT (24-aug) ---
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting of the last header cell (BS1) onto the new header
# cell (BT1) so it keeps the bold/centered/bordered header style.
$ws1.Range("BS1").Copy()
$ws1.Range("BT1").PasteSpecial(-4122)  # xlPasteFormats
$ws1.Range("BT1").Value = "24-aug"

$bt = @{
    2  = 100.5
    3  = 95
    4  = 88.98999999999999
    5  = 85.66
    6  = 83.61
    7  = 86.53
    8  = 90
    9  = 82.51000000000001
    10 = 67.06
    11 = 4.56
    12 = -0.01
    13 = -0.07000000000000001
    14 = -1.98
    15 = -7.45
    16 = -7.4
    17 = -0.1
    18 = 0
    19 = 12.93
    20 = 61.15
    21 = 98.56999999999999
    22 = 114.99
    23 = 114.32
    24 = 107.5
    25 = 100.52
}

foreach ($row in $bt.Keys) {
    $ws1.Cells.Item($row, 72).Value = $bt[$row]
}

# --- Sheet "Gaz": append row 69 ---
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Range("A69").NumberFormat = "@"
$ws2.Range("A69").Value = "2025-08-22"
$ws2.Range("A69").Style = "Normal"
$ws2.Range("B69").Value = 32.2

# --- Sheet "CO2": append row 69 (A69 has date, B69 stays blank) ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A69").NumberFormat = "@"
$ws3.Range("A69").Value = "2025-08-22"
$ws3.Range("A69").Style = "Normal"
